$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------
# 1) "... stdio.h file / in the ..." -> merge the " file" and " in the"
#    runs into a single run, leaving the rest of the paragraph intact.
# ---------------------------------------------------------------------
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*stdio.h file in the directory*") {
        $p1 = $p
        break
    }
}

$xml1 = '<w:p ' + $wNs + ' w14:paraId="41C401EA" w14:textId="00A06652" w:rsidR="0044155A" w:rsidRPr="0044155A" w:rsidRDefault="00F44DF8" w:rsidP="003236BC">' + `
    '<w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">It looks for the </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>stdio.h</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> file in the</w:t></w:r>' + `
    '<w:r w:rsidR="00E954CB"><w:t xml:space="preserve"> directory &#8220;</w:t></w:r>' + `
    '<w:r w:rsidR="00E954CB" w:rsidRPr="00E954CB"><w:t>/</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r w:rsidR="00E954CB" w:rsidRPr="00E954CB"><w:t>usr</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r w:rsidR="00E954CB" w:rsidRPr="00E954CB"><w:t>/include</w:t></w:r>' + `
    '<w:r w:rsidR="00E954CB"><w:t>&#8221;</w:t></w:r>' + `
    '<w:r><w:t>.</w:t></w:r>' + `
    '</w:p>'
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) Question 1.3 answer: "<Your answer here>" -> the real answer about
#    missing enq(double)/deq() declarations in "queue.h".
#    (This is the FIRST "<Your answer here>" paragraph in the doc.)
# ---------------------------------------------------------------------
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<Your answer here>`r") {
        $p2 = $p
        break
    }
}

$xml2 = '<w:p ' + $wNs + ' w14:paraId="01AF7449" w14:textId="77877FF0" w:rsidR="0044155A" w:rsidRPr="0044155A" w:rsidRDefault="0044155A" w:rsidP="0044155A">' + `
    '<w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="0044155A"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Missing function declarations/prototypes for </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>enq</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">(double) and </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>deq</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>()</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#8220;</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>queue.h</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#8221;.</w:t></w:r>' + `
    '</w:p>'
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) Question 1.4 answer: "<Your answer here>" -> two code lines
#    ("void enq(double);" and "double deq();") separated by a blank
#    paragraph, and the paragraph's own center/justify alignment is
#    dropped. (After step 2, this is again the FIRST remaining
#    "<Your answer here>" paragraph, since Q1.3's was just replaced.)
# ---------------------------------------------------------------------
$p3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<Your answer here>`r") {
        $p3 = $p
        break
    }
}

$xml3 = '<w:p ' + $wNs + ' w14:paraId="3E61C2FE" w14:textId="77777777" w:rsidR="0044155A" w:rsidRPr="0044155A" w:rsidRDefault="0044155A" w:rsidP="0044155A">' + `
    '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">void </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>enq</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(double);</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p ' + $wNs + ' w14:paraId="3E61C2FF" w14:textId="77777777" w:rsidR="0044155A" w:rsidRDefault="0044155A" w:rsidP="0044155A"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
    '<w:p ' + $wNs + ' w14:paraId="3E61C300" w14:textId="77777777" w:rsidR="0044155A" w:rsidRDefault="0044155A" w:rsidP="0044155A"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">double </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>deq</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>();</w:t></w:r>' + `
    '</w:p>'
$p3.Range.InsertXML($xml3)

Write-Host "Edits applied."

$cnt = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<Your answer here>`r") {
        $cnt = $cnt + 1
    }
}
Write-Host "Remaining placeholders:" $cnt
